$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("test-meta")

# Insert a new column before column D, shifting existing D..G to E..H
$ws.Range("D1").EntireColumn.Insert()

# Set the header text for the newly inserted column D
$ws.Range("D1").Value = "is_normal_for_donor"

# Set the width of the new column D to match the source file
$ws.Range("D1").EntireColumn.ColumnWidth = 17.46

# Select D2 as the active cell (matches the saved selection in the file)
$ws.Range("D2").Select() | Out-Null
